# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for Cebollín (Femacal de La Calera,
# Coquimbo / Provincia de Quillota) ahead of the existing row 782, shifting
# all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 782..840 down to 783..841, leaving a blank row 782 to populate.
$ws.Rows(782).Insert()

$row = 782
$ws.Cells.Item($row, 1).Value = 3
$ws.Cells.Item($row, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item($row, 3).Value = 'Coquimbo'
$ws.Cells.Item($row, 4).Value = 45223
$ws.Cells.Item($row, 5).Value = 5
$ws.Cells.Item($row, 6).Value = 100112037
$ws.Cells.Item($row, 7).Value = 'Cebollín'
$ws.Cells.Item($row, 8).Value = 'Sin especificar'
$ws.Cells.Item($row, 9).Value = 'Primera'
$ws.Cells.Item($row, 10).Value = 120
$ws.Cells.Item($row, 11).Value = 4000
$ws.Cells.Item($row, 12).Value = 4000
$ws.Cells.Item($row, 13).Value = 4000
$ws.Cells.Item($row, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item($row, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item($row, 16).Value = 111
$ws.Cells.Item($row, 17).Value = 36
$ws.Cells.Item($row, 18).Value = 'Hortaliza'
